# "component approach to pulling URLs and data-anas"
# The "Data Attributes Pages" sheet used to keep a "URL" header in column C
# followed by the two BMO URLs (with hyperlink style) in C3/C4. The new
# layout drops the header and the hyperlinked URLs themselves now live in
# column A (A2/A3), with a couple of extra styled-but-empty rows below
# (A4/A5) reserved for upcoming data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the old hyperlinks that lived on C3/C4 before we move anything around.
$ws.Range("C3").Hyperlinks.Delete()

# Wipe out the old "URL" header + the two old hyperlinked cells in column C.
$ws.Range("C2:C4").Clear()

# Write the two URLs into column A.
$ws.Range("A2").Value2 = "https://www.bmo.com/main/personal"
$ws.Range("A3").Value2 = "https://www.bmo.com/main/personal/credit-cards"

# Re-create the hyperlinks, now anchored on column A.
$ws.Hyperlinks.Add($ws.Range("A2"), "https://www.bmo.com/main/personal")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.bmo.com/main/personal/credit-cards")

# A2:A5 all carry the Hyperlink cell style (A4/A5 are reserved/empty for now).
$ws.Range("A2:A5").Style = "Hyperlink"
